# Updated cryptos list with new price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text (string) representation so
# values like "333.10" are not silently turned into the number 333.1
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.553.57"
$ws.Range("E2").Value = "  +3.82%  "

$ws.Range("D3").Value = "2.765.30"
$ws.Range("E3").Value = "  +4.91%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "116.13"
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").Value = "333.10"
$ws.Range("E6").Value = "  +2.79%  "

$ws.Range("D7").Value = "0.539"
$ws.Range("E7").Value = "  +2.12%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +5.38%  "

$ws.Range("D10").Value = "41.87"
$ws.Range("E10").Value = "  +5.07%  "

$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  +6.89%  "

$ws.Range("D12").Value = "20.26"
$ws.Range("E12").Value = "  +2.24%  "

$ws.Range("D13").Value = "0.130"
$ws.Range("E13").Value = "  +2.33%  "

$ws.Range("D14").Value = "7.66"
$ws.Range("E14").Value = "  +4.67%  "

$ws.Range("D15").Value = "3.193.10"
$ws.Range("E15").Value = "  +4.49%  "

$ws.Range("D16").Value = "2.773.69"
$ws.Range("E16").Value = "  +4.99%  "

$ws.Range("D17").Value = "0.891"
$ws.Range("E17").Value = "  +3.55%  "

$ws.Range("D18").Value = "51.559.00"
$ws.Range("E18").Value = "  +3.96%  "

$ws.Range("D19").Value = "3.31"
$ws.Range("E19").Value = "  +11.11%  "

$ws.Range("D20").Value = "13.51"
$ws.Range("E20").Value = "  +4.52%  "

$ws.Range("E21").Value = "  +2.26%  "

$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +3.27%  "

$ws.Range("D23").Value = "278.40"
$ws.Range("E23").Value = "  +3.04%  "

$ws.Range("D24").Value = "69.74"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("E25").Value = "  +6.06%  "

$ws.Range("D26").Value = "26.84"
$ws.Range("E26").Value = "  +1.94%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").Value = "35.14"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").Value = "50.14"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("D33").Value = "5.57"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").Value = "19.03"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("D38").Value = "2.09"
$ws.Range("E38").Value = "  +2.05%  "

$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("D40").Value = "0.0353"
$ws.Range("E40").Value = "  +8.63%  "

$ws.Range("D41").Value = "127.41"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("E44").Value = "  +6.87%  "

$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +16.06%  "

$ws.Range("D46").Value = "2.090.38"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +2.76%  "

$ws.Range("E48").Value = "  +3.27%  "

$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +5.95%  "

$ws.Range("D50").Value = "9.00"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("D51").Value = "60.28"
$ws.Range("E51").Value = "  +1.88%  "

# Rows 42/43 swapped position (Stellar now ranked above EnergySwap)
# Row 42: EnergySwap -> Stellar
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  +2.99%  "

# Row 43: Stellar -> EnergySwap
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "23.04"
$ws.Range("E43").Value = "  +3.12%  "
